$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Insert two new blank columns at D (shifts old D:K -> F:M) ---
$ws.Range("D1:E1").EntireColumn.Insert()

# --- Step 2: Copy number formats into the newly inserted D:E columns from F (old D) ---
$ws.Range("F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F8:F102").Copy()
$ws.Range("D8:E102").PasteSpecial(-4122)

# --- Step 3: Append two new columns L:M, copying formats from K ---
$ws.Range("K7").Copy()
$ws.Range("L7:M7").PasteSpecial(-4122)
$ws.Range("K8:K102").Copy()
$ws.Range("L8:M102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 4: Write final values for every data cell D:M across all data rows ---
$ws.Range("D7").Value = 43464
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43282
$ws.Range("G7").Value = 43191
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43009
$ws.Range("J7").Value = 42918
$ws.Range("K7").Value = 42827
$ws.Range("L7").Value = 42736
$ws.Range("M7").Value = 42645
$ws.Range("D8").Value = 5207300
$ws.Range("E8").Value = 1329100
$ws.Range("F8").Value = 1321700
$ws.Range("G8").Value = 1287200
$ws.Range("H8").Value = 4664600
$ws.Range("I8").Value = 1206100
$ws.Range("J8").Value = 1184000
$ws.Range("K8").Value = 1130600
$ws.Range("L8").Value = 985700
$ws.Range("M8").Value = 1035800
$ws.Range("D9").Value = 3459900
$ws.Range("E9").Value = 946700
$ws.Range("F9").Value = 941300
$ws.Range("G9").Value = 900100
$ws.Range("H9").Value = 3917100
$ws.Range("I9").Value = 859700
$ws.Range("J9").Value = 842000
$ws.Range("K9").Value = 793400
$ws.Range("L9").Value = 707500
$ws.Range("M9").Value = 744300
$ws.Range("D10").Value = 1747400
$ws.Range("E10").Value = 382400
$ws.Range("F10").Value = 380400
$ws.Range("G10").Value = 387100
$ws.Range("H10").Value = 747500
$ws.Range("I10").Value = 346500
$ws.Range("J10").Value = 342000
$ws.Range("K10").Value = 337200
$ws.Range("L10").Value = 278200
$ws.Range("M10").Value = 291500
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = ""
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = "NA"
$ws.Range("M12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = "NA"
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("D15").Value = 108000
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "NA"
$ws.Range("G15").Value = "NA"
$ws.Range("H15").Value = 94200
$ws.Range("I15").Value = "NA"
$ws.Range("J15").Value = "NA"
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = ""
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = ""
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = ""
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = ""
$ws.Range("D17").Value = 4984400
$ws.Range("E17").Value = 1276300
$ws.Range("F17").Value = 1260000
$ws.Range("G17").Value = 1207500
$ws.Range("H17").Value = 4438500
$ws.Range("I17").Value = 1153100
$ws.Range("J17").Value = 1120500
$ws.Range("K17").Value = 1058100
$ws.Range("L17").Value = 955500
$ws.Range("M17").Value = 994400
$ws.Range("D18").Value = 222900
$ws.Range("E18").Value = 52800
$ws.Range("F18").Value = 61700
$ws.Range("G18").Value = 79700
$ws.Range("H18").Value = 226100
$ws.Range("I18").Value = 53000
$ws.Range("J18").Value = 63500
$ws.Range("K18").Value = 72500
$ws.Range("L18").Value = 30200
$ws.Range("M18").Value = 41400
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = ""
$ws.Range("I19").Value = ""
$ws.Range("J19").Value = ""
$ws.Range("K19").Value = ""
$ws.Range("L19").Value = ""
$ws.Range("M19").Value = ""
$ws.Range("D20").Value = 400
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 100
$ws.Range("G20").Value = 200
$ws.Range("H20").Value = 600
$ws.Range("I20").Value = 100
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 100
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = 100
$ws.Range("D21").Value = 334000
$ws.Range("E21").Value = 80900
$ws.Range("F21").Value = 88800
$ws.Range("G21").Value = 106700
$ws.Range("H21").Value = 323400
$ws.Range("I21").Value = 77900
$ws.Range("J21").Value = 87100
$ws.Range("K21").Value = 95200
$ws.Range("L21").Value = 50700
$ws.Range("M21").Value = 62800
$ws.Range("D22").Value = 27400
$ws.Range("E22").Value = 7400
$ws.Range("F22").Value = 6500
$ws.Range("G22").Value = 6100
$ws.Range("H22").Value = 21200
$ws.Range("I22").Value = 5600
$ws.Range("J22").Value = 5100
$ws.Range("K22").Value = 4700
$ws.Range("L22").Value = 3800
$ws.Range("M22").Value = 3700
$ws.Range("D23").Value = 195800
$ws.Range("E23").Value = 45400
$ws.Range("F23").Value = 55300
$ws.Range("G23").Value = 73800
$ws.Range("H23").Value = 205500
$ws.Range("I23").Value = 47600
$ws.Range("J23").Value = 58500
$ws.Range("K23").Value = 67900
$ws.Range("L23").Value = 26500
$ws.Range("M23").Value = 37900
$ws.Range("D24").Value = 37300
$ws.Range("E24").Value = 7900
$ws.Range("F24").Value = 13600
$ws.Range("G24").Value = 7200
$ws.Range("H24").Value = 65800
$ws.Range("I24").Value = 16100
$ws.Range("J24").Value = 17500
$ws.Range("K24").Value = 21600
$ws.Range("L24").Value = 9500
$ws.Range("M24").Value = 14000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = 158500
$ws.Range("E26").Value = 37500
$ws.Range("F26").Value = 41700
$ws.Range("G26").Value = 66600
$ws.Range("H26").Value = 139700
$ws.Range("I26").Value = 31500
$ws.Range("J26").Value = 41000
$ws.Range("K26").Value = 46300
$ws.Range("L26").Value = 17000
$ws.Range("M26").Value = 23900
$ws.Range("D27").Value = 158500
$ws.Range("E27").Value = 37500
$ws.Range("F27").Value = 41700
$ws.Range("G27").Value = 66600
$ws.Range("H27").Value = 139700
$ws.Range("I27").Value = 31500
$ws.Range("J27").Value = 41000
$ws.Range("K27").Value = 46300
$ws.Range("L27").Value = 17000
$ws.Range("M27").Value = 23900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = 18700
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
$ws.Range("K29").Value = "NA"
$ws.Range("L29").Value = "NA"
$ws.Range("M29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("D32").Value = -400
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = -100
$ws.Range("G32").Value = -200
$ws.Range("H32").Value = -600
$ws.Range("I32").Value = -100
$ws.Range("J32").Value = -100
$ws.Range("K32").Value = -100
$ws.Range("L32").Value = -100
$ws.Range("M32").Value = -100
$ws.Range("D33").Value = 158500
$ws.Range("E33").Value = 37500
$ws.Range("F33").Value = 41700
$ws.Range("G33").Value = 66600
$ws.Range("H33").Value = 158400
$ws.Range("I33").Value = 31500
$ws.Range("J33").Value = 41000
$ws.Range("K33").Value = 46300
$ws.Range("L33").Value = 17000
$ws.Range("M33").Value = 23900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("D35").Value = 158500
$ws.Range("E35").Value = 37500
$ws.Range("F35").Value = 41700
$ws.Range("G35").Value = 66600
$ws.Range("H35").Value = 158400
$ws.Range("I35").Value = 31500
$ws.Range("J35").Value = 41000
$ws.Range("K35").Value = 46300
$ws.Range("L35").Value = 17000
$ws.Range("M35").Value = 23900
$ws.Range("D38").Value = 43464
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43282
$ws.Range("G38").Value = 43191
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43009
$ws.Range("J38").Value = 42918
$ws.Range("K38").Value = 42827
$ws.Range("L38").Value = 42736
$ws.Range("M38").Value = 42645
$ws.Range("D39").Value = ""
$ws.Range("E39").Value = ""
$ws.Range("F39").Value = ""
$ws.Range("G39").Value = ""
$ws.Range("H39").Value = ""
$ws.Range("I39").Value = ""
$ws.Range("J39").Value = ""
$ws.Range("K39").Value = ""
$ws.Range("L39").Value = ""
$ws.Range("M39").Value = ""
$ws.Range("D40").Value = ""
$ws.Range("E40").Value = ""
$ws.Range("F40").Value = ""
$ws.Range("G40").Value = ""
$ws.Range("H40").Value = ""
$ws.Range("I40").Value = ""
$ws.Range("J40").Value = ""
$ws.Range("K40").Value = ""
$ws.Range("L40").Value = ""
$ws.Range("M40").Value = ""
$ws.Range("D41").Value = 1600
$ws.Range("E41").Value = 17200
$ws.Range("F41").Value = 23200
$ws.Range("G41").Value = 21400
$ws.Range("H41").Value = 19500
$ws.Range("I41").Value = 18900
$ws.Range("J41").Value = 14800
$ws.Range("K41").Value = 21300
$ws.Range("L41").Value = 12500
$ws.Range("M41").Value = 50300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("D43").Value = 40600
$ws.Range("E43").Value = 34600
$ws.Range("F43").Value = 27600
$ws.Range("G43").Value = 26700
$ws.Range("H43").Value = 25900
$ws.Range("I43").Value = 23200
$ws.Range("J43").Value = 19800
$ws.Range("K43").Value = 24700
$ws.Range("L43").Value = 25200
$ws.Range("M43").Value = 21800
$ws.Range("D44").Value = 264400
$ws.Range("E44").Value = 253000
$ws.Range("F44").Value = 250600
$ws.Range("G44").Value = 239600
$ws.Range("H44").Value = 229500
$ws.Range("I44").Value = 222200
$ws.Range("J44").Value = 218400
$ws.Range("K44").Value = 204300
$ws.Range("L44").Value = 204500
$ws.Range("M44").Value = 195200
$ws.Range("D45").Value = 27300
$ws.Range("E45").Value = 37500
$ws.Range("F45").Value = 28700
$ws.Range("G45").Value = 27300
$ws.Range("H45").Value = 24600
$ws.Range("I45").Value = 25600
$ws.Range("J45").Value = 23400
$ws.Range("K45").Value = 23200
$ws.Range("L45").Value = 21900
$ws.Range("M45").Value = 24500
$ws.Range("D46").Value = 333800
$ws.Range("E46").Value = 342400
$ws.Range("F46").Value = 330000
$ws.Range("G46").Value = 315000
$ws.Range("H46").Value = 299500
$ws.Range("I46").Value = 289900
$ws.Range("J46").Value = 276400
$ws.Range("K46").Value = 273500
$ws.Range("L46").Value = 264000
$ws.Range("M46").Value = 291800
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("D48").Value = 766400
$ws.Range("E48").Value = 773300
$ws.Range("F48").Value = 776300
$ws.Range("G48").Value = 738700
$ws.Range("H48").Value = 713000
$ws.Range("I48").Value = 690800
$ws.Range("J48").Value = 681100
$ws.Range("K48").Value = 653400
$ws.Range("L48").Value = 604700
$ws.Range("M48").Value = 577400
$ws.Range("D49").Value = 562900
$ws.Range("E49").Value = 563200
$ws.Range("F49").Value = 563600
$ws.Range("G49").Value = 563900
$ws.Range("H49").Value = 564300
$ws.Range("I49").Value = 564600
$ws.Range("J49").Value = 565000
$ws.Range("K49").Value = 565300
$ws.Range("L49").Value = 565700
$ws.Range("M49").Value = 566000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("D52").Value = 12500
$ws.Range("E52").Value = 16000
$ws.Range("F52").Value = 13900
$ws.Range("G52").Value = 9500
$ws.Range("H52").Value = 4800
$ws.Range("I52").Value = 5900
$ws.Range("J52").Value = 7800
$ws.Range("K52").Value = 5400
$ws.Range("L52").Value = 5500
$ws.Range("M52").Value = 20100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("D54").Value = 1675600
$ws.Range("E54").Value = 1695000
$ws.Range("F54").Value = 1683800
$ws.Range("G54").Value = 1627100
$ws.Range("H54").Value = 1581600
$ws.Range("I54").Value = 1551200
$ws.Range("J54").Value = 1530200
$ws.Range("K54").Value = 1497600
$ws.Range("L54").Value = 1439900
$ws.Range("M54").Value = 1455400
$ws.Range("D55").Value = ""
$ws.Range("E55").Value = ""
$ws.Range("F55").Value = ""
$ws.Range("G55").Value = ""
$ws.Range("H55").Value = ""
$ws.Range("I55").Value = ""
$ws.Range("J55").Value = ""
$ws.Range("K55").Value = ""
$ws.Range("L55").Value = ""
$ws.Range("M55").Value = ""
$ws.Range("D56").Value = ""
$ws.Range("E56").Value = ""
$ws.Range("F56").Value = ""
$ws.Range("G56").Value = ""
$ws.Range("H56").Value = ""
$ws.Range("I56").Value = ""
$ws.Range("J56").Value = ""
$ws.Range("K56").Value = ""
$ws.Range("L56").Value = ""
$ws.Range("M56").Value = ""
$ws.Range("D57").Value = 254000
$ws.Range("E57").Value = 245900
$ws.Range("F57").Value = 265600
$ws.Range("G57").Value = 270400
$ws.Range("H57").Value = 244900
$ws.Range("I57").Value = 179500
$ws.Range("J57").Value = 181900
$ws.Range("K57").Value = 191800
$ws.Range("L57").Value = 157600
$ws.Range("M57").Value = 162300
$ws.Range("D58").Value = 7400
$ws.Range("E58").Value = 7400
$ws.Range("F58").Value = 14400
$ws.Range("G58").Value = 12000
$ws.Range("H58").Value = 9200
$ws.Range("I58").Value = 8800
$ws.Range("J58").Value = 14900
$ws.Range("K58").Value = 10600
$ws.Range("L58").Value = 12400
$ws.Range("M58").Value = 9400
$ws.Range("D59").Value = 48600
$ws.Range("E59").Value = 43400
$ws.Range("F59").Value = 38400
$ws.Range("G59").Value = 33600
$ws.Range("H59").Value = 45600
$ws.Range("I59").Value = 101900
$ws.Range("J59").Value = 85100
$ws.Range("K59").Value = 97700
$ws.Range("L59").Value = 89200
$ws.Range("M59").Value = 73300
$ws.Range("D60").Value = 310000
$ws.Range("E60").Value = 296800
$ws.Range("F60").Value = 318400
$ws.Range("G60").Value = 315900
$ws.Range("H60").Value = 299700
$ws.Range("I60").Value = 290200
$ws.Range("J60").Value = 281800
$ws.Range("K60").Value = 300100
$ws.Range("L60").Value = 259200
$ws.Range("M60").Value = 245000
$ws.Range("D61").Value = 572600
$ws.Range("E61").Value = 555700
$ws.Range("F61").Value = 579600
$ws.Range("G61").Value = 490800
$ws.Range("H61").Value = 473500
$ws.Range("I61").Value = 475800
$ws.Range("J61").Value = 434700
$ws.Range("K61").Value = 404400
$ws.Range("L61").Value = 372400
$ws.Range("M61").Value = 320400
$ws.Range("D62").Value = 203800
$ws.Range("E62").Value = 201200
$ws.Range("F62").Value = 184600
$ws.Range("G62").Value = 171900
$ws.Range("H62").Value = 157700
$ws.Range("I62").Value = 168600
$ws.Range("J62").Value = 161300
$ws.Range("K62").Value = 149200
$ws.Range("L62").Value = 135500
$ws.Range("M62").Value = 130600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("D66").Value = 1086400
$ws.Range("E66").Value = 1053700
$ws.Range("F66").Value = 1082600
$ws.Range("G66").Value = 978700
$ws.Range("H66").Value = 930900
$ws.Range("I66").Value = 934600
$ws.Range("J66").Value = 877900
$ws.Range("K66").Value = 853700
$ws.Range("L66").Value = 767000
$ws.Range("M66").Value = 696000
$ws.Range("D67").Value = ""
$ws.Range("E67").Value = ""
$ws.Range("F67").Value = ""
$ws.Range("G67").Value = ""
$ws.Range("H67").Value = ""
$ws.Range("I67").Value = ""
$ws.Range("J67").Value = ""
$ws.Range("K67").Value = ""
$ws.Range("L67").Value = ""
$ws.Range("M67").Value = ""
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("D72").Value = -69200
$ws.Range("E72").Value = -16900
$ws.Range("F72").Value = -39100
$ws.Range("G72").Value = 14200
$ws.Range("H72").Value = 30600
$ws.Range("I72").Value = 2300
$ws.Range("J72").Value = 42800
$ws.Range("K72").Value = 41800
$ws.Range("L72").Value = 75500
$ws.Range("M72").Value = 164900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("D76").Value = 589200
$ws.Range("E76").Value = 641300
$ws.Range("F76").Value = 601200
$ws.Range("G76").Value = 648500
$ws.Range("H76").Value = 650700
$ws.Range("I76").Value = 616600
$ws.Range("J76").Value = 652400
$ws.Range("K76").Value = 643900
$ws.Range("L76").Value = 672900
$ws.Range("M76").Value = 759300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("D80").Value = 43464
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43282
$ws.Range("G80").Value = 43191
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43009
$ws.Range("J80").Value = 42918
$ws.Range("K80").Value = 42827
$ws.Range("L80").Value = 42736
$ws.Range("M80").Value = 42645
$ws.Range("D81").Value = 158500
$ws.Range("E81").Value = 37500
$ws.Range("F81").Value = 41700
$ws.Range("G81").Value = 66600
$ws.Range("H81").Value = 158400
$ws.Range("I81").Value = 31500
$ws.Range("J81").Value = 41000
$ws.Range("K81").Value = 46300
$ws.Range("L81").Value = 17000
$ws.Range("M81").Value = 23900
$ws.Range("D82").Value = ""
$ws.Range("E82").Value = ""
$ws.Range("F82").Value = ""
$ws.Range("G82").Value = ""
$ws.Range("H82").Value = ""
$ws.Range("I82").Value = ""
$ws.Range("J82").Value = ""
$ws.Range("K82").Value = ""
$ws.Range("L82").Value = ""
$ws.Range("M82").Value = ""
$ws.Range("D83").Value = 110700
$ws.Range("E83").Value = 28100
$ws.Range("F83").Value = 27000
$ws.Range("G83").Value = 26800
$ws.Range("H83").Value = 96700
$ws.Range("I83").Value = 24700
$ws.Range("J83").Value = 23500
$ws.Range("K83").Value = 22600
$ws.Range("L83").Value = 20400
$ws.Range("M83").Value = 21200
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("D89").Value = 294400
$ws.Range("E89").Value = 64400
$ws.Range("F89").Value = 66600
$ws.Range("G89").Value = 104500
$ws.Range("H89").Value = 309600
$ws.Range("I89").Value = 79100
$ws.Range("J89").Value = 65300
$ws.Range("K89").Value = 114500
$ws.Range("L89").Value = 58300
$ws.Range("M89").Value = 48100
$ws.Range("D90").Value = ""
$ws.Range("E90").Value = ""
$ws.Range("F90").Value = ""
$ws.Range("G90").Value = ""
$ws.Range("H90").Value = ""
$ws.Range("I90").Value = ""
$ws.Range("J90").Value = ""
$ws.Range("K90").Value = ""
$ws.Range("L90").Value = ""
$ws.Range("M90").Value = ""
$ws.Range("D91").Value = -177100
$ws.Range("E91").Value = -44500
$ws.Range("F91").Value = -59800
$ws.Range("G91").Value = -44200
$ws.Range("H91").Value = -198600
$ws.Range("I91").Value = -40700
$ws.Range("J91").Value = -60500
$ws.Range("K91").Value = -57200
$ws.Range("L91").Value = -38400
$ws.Range("M91").Value = -57500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("D94").Value = -177100
$ws.Range("E94").Value = -44500
$ws.Range("F94").Value = -59800
$ws.Range("G94").Value = -44200
$ws.Range("H94").Value = -198600
$ws.Range("I94").Value = -40700
$ws.Range("J94").Value = -60500
$ws.Range("K94").Value = -57200
$ws.Range("L94").Value = -38400
$ws.Range("M94").Value = -57500
$ws.Range("D95").Value = ""
$ws.Range("E95").Value = ""
$ws.Range("F95").Value = ""
$ws.Range("G95").Value = ""
$ws.Range("H95").Value = ""
$ws.Range("I95").Value = ""
$ws.Range("J95").Value = ""
$ws.Range("K95").Value = ""
$ws.Range("L95").Value = ""
$ws.Range("M95").Value = ""
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("D100").Value = -134500
$ws.Range("E100").Value = -25200
$ws.Range("F100").Value = -5000
$ws.Range("G100").Value = -58500
$ws.Range("H100").Value = -104000
$ws.Range("I100").Value = -34300
$ws.Range("J100").Value = -11300
$ws.Range("K100").Value = -48500
$ws.Range("L100").Value = -57700
$ws.Range("M100").Value = -18800
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("D102").Value = -17200
$ws.Range("E102").Value = -5300
$ws.Range("F102").Value = 1800
$ws.Range("G102").Value = 1900
$ws.Range("H102").Value = 7000
$ws.Range("I102").Value = 4100
$ws.Range("J102").Value = -6500
$ws.Range("K102").Value = 8900
$ws.Range("L102").Value = -37800
$ws.Range("M102").Value = -28200

Write-Host "Done writing values"